{"js": "// The template contains a Word field (fldChar begin / instrText \" m:null.setConserveRatio(true) \" /\n// fldChar end) that encodes an M2Doc expression. This edit rewrites it to the equivalent\n// literal-text bracket notation \"{m:null.setConserveRatio(true)}\" used by\n// TokenIteratorFieldRewriterSplit, split across w:r/w:t runs (the \"null\" and\n// \".setConserveRatio(true)\" runs keep their orange run formatting).\n\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\n// Locate the field whose code is the M2Doc \"m:null.setConserveRatio(true)\" expression.\nfor (let i = 0; i < fields.items.length; i++) {\n  fields.items[i].load(\"code\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < fields.items.length; i++) {\n  const f = fields.items[i];\n  if (f.code && f.code.indexOf(\"setConserveRatio\") !== -1) {\n    target = f;\n    break;\n  }\n}\nif (!target && fields.items.length > 0) {\n  target = fields.items[0];\n}\n\n// Grab the paragraph that currently hosts the field (the field's code/result text is not part\n// of paragraph.text, so find it by elimination: it is the empty paragraph immediately before\n// the \"End of demonstration.\" paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nlet hostParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"\" && i + 1 < paragraphs.items.length &&\n      paragraphs.items[i + 1].text.indexOf(\"End of demonstration\") !== -1) {\n    hostParagraph = p;\n    break;\n  }\n}\n\n// Delete the field (removes the fldChar/instrText runs) while preserving the paragraph mark\n// itself (identity/formatting of the <w:p> stay untouched).\ntarget.delete();\nawait context.sync();\n\n// Insert the literal-text replacement runs into the now-empty paragraph:\n// \"{\" + \"m\" + \":\" + \"null\" + \".setConserveRatio(true)\" + \"}\" -- the last two runs keep the\n// original orange accent-color run formatting.\nconst contentRange = hostParagraph.getRange(\"Content\");\n\nconst runsOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.setConserveRatio(true)</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ncontentRange.insertOoxml(runsOoxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# The template contains a Word field (fldChar begin / instrText \" m:null.setConserveRatio(true) \" /\n# fldChar end) that encodes an M2Doc expression. This edit rewrites it to the equivalent\n# literal-text bracket notation \"{m:null.setConserveRatio(true)}\" used by\n# TokenIteratorFieldRewriterSplit, split across w:r/w:t runs (the \"null\" and\n# \".setConserveRatio(true)\" runs keep their orange run formatting).\n\n$d = $word.ActiveDocument\n\n# Locate the field that holds the M2Doc \"m:null.setConserveRatio(true)\" expression.\n$target = $null\nforeach ($f in $d.Fields) {\n  if ($f.Code.Text -match \"setConserveRatio\") {\n    $target = $f\n    break\n  }\n}\nif ($target -eq $null -and $d.Fields.Count -gt 0) {\n  $target = $d.Fields.Item(1)\n}\n\n# Find the paragraph hosting this field by matching the field code's start offset against each\n# paragraph's range (the field's code/result text isn't part of paragraph.Range.Text, so we\n# can't find it by text content).\n$hostPara = $null\nforeach ($p in $d.Paragraphs) {\n  if ($target.Code.Start -ge $p.Range.Start -and $target.Code.Start -lt $p.Range.End) {\n    $hostPara = $p\n    break\n  }\n}\n\n# Capture the paragraph's own opening tag (its identity / rsid / paragraph-mark formatting) so\n# the rewritten paragraph keeps it unchanged.\n$openTag = '<w:p>'\n$hostXml = $hostPara.Range.WordOpenXML\nif ($hostXml -match '<w:p( [^>]*)?>') {\n  $openTag = $matches[0]\n}\n\n$hostRange = $hostPara.Range\n\n# Delete the field (removes the fldChar/instrText runs) while leaving the paragraph mark itself\n# untouched.\n$target.Delete()\n\n# Build the literal-text replacement runs: \"{\" + \"m\" + \":\" + \"null\" + \".setConserveRatio(true)\" + \"}\"\n# -- the last two runs keep the original orange accent-color character formatting.\n$newParaXml = $openTag +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>null</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.setConserveRatio(true)</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>'\n\n$fullXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' + $newParaXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n$hostRange.InsertXML($fullXml)\n"}
